$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Add the new "Open" column header in J1, matching the style of the other
# required (grey) columns such as A1.
$ws.Range("J1").Value = "Open"
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Document the new column with a cell comment.
$ws.Range("J1").AddComment("Set this to 0 to import the ticket as closed, 1 to import as open.")

# Match the updated selection recorded for the sheet.
[void]$ws.Range("I15").Select()
